$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 13.46456841828657
$ws.Range("D2").Value = 7.237647162246525
$ws.Range("E2").Value = 14.00369035395977
$ws.Range("F2").Value = 47.40230617506971
$ws.Range("G2").Value = 3.711498604868209
$ws.Range("J2").Value = 11.29078040900452
$ws.Range("L2").Value = 9.595067975462683
$ws.Range("O2").Value = 36.86143484149457

$ws.Range("C3").Value = 13.44089697468714
$ws.Range("D3").Value = 7.249130751239226
$ws.Range("E3").Value = 14.02303258333953
$ws.Range("F3").Value = 47.12862184932497
$ws.Range("G3").Value = 3.715841198117704
$ws.Range("J3").Value = 11.31216645505454
$ws.Range("L3").Value = 9.615335707911301
$ws.Range("O3").Value = 36.66990447118242

$ws.Range("C4").Value = 13.42938265462392
$ws.Range("D4").Value = 7.256736956643145
$ws.Range("E4").Value = 14.03716565015119
$ws.Range("F4").Value = 46.97262192089115
$ws.Range("G4").Value = 3.718646908430662
$ws.Range("J4").Value = 11.3269109368498
$ws.Range("L4").Value = 9.628643202278969
$ws.Range("O4").Value = 36.56157404316798

$ws.Range("C5").Value = 13.42545186324616
$ws.Range("D5").Value = 7.259976432875158
$ws.Range("E5").Value = 14.04349205974692
$ws.Range("F5").Value = 46.91212518869642
$ws.Range("G5").Value = 3.719825436473871
$ws.Range("J5").Value = 11.33332499675213
$ws.Range("L5").Value = 9.634283630440986
$ws.Range("O5").Value = 36.51978946809265

$ws.Range("C6").Value = 13.42484518937866
$ws.Range("D6").Value = 7.260522802078719
$ws.Range("E6").Value = 14.04457679044524
$ws.Range("F6").Value = 46.90226669068198
$ws.Range("G6").Value = 3.720023258950125
$ws.Range("J6").Value = 11.33441453733435
$ws.Range("L6").Value = 9.635233371458447
$ws.Range("O6").Value = 36.51299456900799

$ws.Range("C7").Value = 13.4293265576303
$ws.Range("D7").Value = 7.25678007859456
$ws.Range("E7").Value = 14.03724867492402
$ws.Range("F7").Value = 46.97179353550013
$ws.Range("G7").Value = 3.718662659852911
$ws.Range("J7").Value = 11.32699579719404
$ws.Range("L7").Value = 9.628718389690167
$ws.Range("O7").Value = 36.56100092582786

$ws.Range("C8").Value = 13.45578106520064
$ws.Range("D8").Value = 7.241491624922933
$ws.Range("E8").Value = 14.00989090444787
$ws.Range("F8").Value = 47.30546942562894
$ws.Range("G8").Value = 3.712967095336456
$ws.Range("J8").Value = 11.29781939222442
$ws.Range("L8").Value = 9.601877470073202
$ws.Range("O8").Value = 36.79349024983676

$ws.Range("C9").Value = 13.53150781846155
$ws.Range("D9").Value = 7.21590476987804
$ws.Range("E9").Value = 13.97416704116746
$ws.Range("F9").Value = 48.0531048756793
$ws.Range("G9").Value = 3.702897338804731
$ws.Range("J9").Value = 11.25341156141796
$ws.Range("L9").Value = 9.556067829806308
$ws.Range("O9").Value = 37.32147138025219

$ws.Range("C10").Value = 13.6014938719128
$ws.Range("D10").Value = 7.19976875209966
$ws.Range("E10").Value = 13.95886894994041
$ws.Range("F10").Value = 48.65592705651143
$ws.Range("G10").Value = 3.696160297802362
$ws.Range("J10").Value = 11.22860026253204
$ws.Range("L10").Value = 9.526541753109228
$ws.Range("O10").Value = 37.7511012954394

$ws.Range("C11").Value = 13.63639464457943
$ws.Range("D11").Value = 7.193002898745478
$ws.Range("E11").Value = 13.95428940830532
$ws.Range("F11").Value = 48.94099179011069
$ws.Range("G11").Value = 3.693237113819939
$ws.Range("J11").Value = 11.21901114985903
$ws.Range("L11").Value = 9.513999979181319
$ws.Range("O11").Value = 37.95507418711851

$ws.Range("C12").Value = 13.65004563244623
$ws.Range("D12").Value = 7.190523199770267
$ws.Range("E12").Value = 13.95289744902856
$ws.Range("F12").Value = 49.05042410277716
$ws.Range("G12").Value = 3.692150384744062
$ws.Range("J12").Value = 11.21562416779595
$ws.Range("L12").Value = 9.50937819002263
$ws.Range("O12").Value = 38.03349089720409

$ws.Range("C13").Value = 13.6470864069972
$ws.Range("D13").Value = 7.191053586949215
$ws.Range("E13").Value = 13.95318201306632
$ws.Range("F13").Value = 49.02679108287165
$ws.Range("G13").Value = 3.692383534264617
$ws.Range("J13").Value = 11.21634275281282
$ws.Range("L13").Value = 9.510367911024186
$ws.Range("O13").Value = 38.01655095156386

$ws.Range("C14").Value = 13.63750904012132
$ws.Range("D14").Value = 7.192797242650474
$ws.Range("E14").Value = 13.95416803274465
$ws.Range("F14").Value = 48.9499655354399
$ws.Range("G14").Value = 3.693147303550291
$ws.Range("J14").Value = 11.2187276061686
$ws.Range("L14").Value = 9.513617188860877
$ws.Range("O14").Value = 37.96150229263643

$ws.Range("C15").Value = 13.63169907231199
$ws.Range("D15").Value = 7.193876004016209
$ws.Range("E15").Value = 13.95481656350215
$ws.Range("F15").Value = 48.90309874599466
$ws.Range("G15").Value = 3.693617763187826
$ws.Range("J15").Value = 11.22022020260248
$ws.Range("L15").Value = 9.515624057686123
$ws.Range("O15").Value = 37.92793509483942

$ws.Range("C16").Value = 13.59927419815431
$ws.Range("D16").Value = 7.200222457274363
$ws.Range("E16").Value = 13.95921612271147
$ws.Range("F16").Value = 48.6375092368538
$ws.Range("G16").Value = 3.696354171179732
$ws.Range("J16").Value = 11.22926110294033
$ws.Range("L16").Value = 9.527379254438779
$ws.Range("O16").Value = 37.73793862966051

$ws.Range("C17").Value = 13.5801634005786
$ws.Range("D17").Value = 7.204262780214714
$ws.Range("E17").Value = 13.96252464355484
$ws.Range("F17").Value = 48.47730286586638
$ws.Range("G17").Value = 3.698069022148583
$ws.Range("J17").Value = 11.23524228606841
$ws.Range("L17").Value = 9.534818266726166
$ws.Range("O17").Value = 37.62353220979392

$ws.Range("C18").Value = 13.56946013472514
$ws.Range("D18").Value = 7.206640755558605
$ws.Range("E18").Value = 13.96465160525664
$ws.Range("F18").Value = 48.38618115253882
$ws.Range("G18").Value = 3.699068688541295
$ws.Range("J18").Value = 11.23884229425479
$ws.Range("L18").Value = 9.539180766894257
$ws.Range("O18").Value = 37.55853491760426

$ws.Range("C19").Value = 13.56588595247155
$ws.Range("D19").Value = 7.207455193937511
$ws.Range("E19").Value = 13.96541022556529
$ws.Range("F19").Value = 48.35550711223614
$ws.Range("G19").Value = 3.699409451859861
$ws.Range("J19").Value = 11.24008863616199
$ws.Range("L19").Value = 9.540672236809339
$ws.Range("O19").Value = 37.53666791909161

$ws.Range("C20").Value = 13.58216793176913
$ws.Range("D20").Value = 7.203827084665054
$ws.Range("E20").Value = 13.96214926343788
$ws.Range("F20").Value = 48.4942515757616
$ws.Range("G20").Value = 3.697885094761881
$ws.Range("J20").Value = 11.23458904057303
$ws.Range("L20").Value = 9.534017703779313
$ws.Range("O20").Value = 37.63562788101444

$ws.Range("C21").Value = 13.64031039274428
$ws.Range("D21").Value = 7.19228285467719
$ws.Range("E21").Value = 13.95386912761617
$ws.Range("F21").Value = 48.97249137689697
$ws.Range("G21").Value = 3.692922418376558
$ws.Range("J21").Value = 11.21802048871059
$ws.Range("L21").Value = 9.512659340151561
$ws.Range("O21").Value = 37.97763989797344

$ws.Range("C22").Value = 13.68084124487542
$ws.Range("D22").Value = 7.18521811846803
$ws.Range("E22").Value = 13.95045219905905
$ws.Range("F22").Value = 49.29366513159245
$ws.Range("G22").Value = 3.68979680195909
$ws.Range("J22").Value = 11.20861542338142
$ws.Range("L22").Value = 9.499443429866568
$ws.Range("O22").Value = 38.20799774611486

$ws.Range("C23").Value = 13.65897959643821
$ws.Range("D23").Value = 7.188944848577767
$ws.Range("E23").Value = 13.95209339001118
$ws.Range("F23").Value = 49.12148544161735
$ws.Range("G23").Value = 3.691454269784962
$ws.Range("J23").Value = 11.21350481907344
$ws.Range("L23").Value = 9.506429168724976
$ws.Range("O23").Value = 38.08444339419753

$ws.Range("C24").Value = 13.58126079885888
$ws.Range("D24").Value = 7.204023890857139
$ws.Range("E24").Value = 13.96231827243731
$ws.Range("F24").Value = 48.48658599641973
$ws.Range("G24").Value = 3.697968205414837
$ws.Range("J24").Value = 11.23488387031035
$ws.Range("L24").Value = 9.534379371277483
$ws.Range("O24").Value = 37.63015700569962

$ws.Range("C25").Value = 13.50848548590462
$ws.Range("D25").Value = 7.222357963614956
$ws.Range("E25").Value = 13.98190979126682
$ws.Range("F25").Value = 47.84118917545359
$ws.Range("G25").Value = 3.70550471752001
$ws.Range("J25").Value = 11.26405308227013
$ws.Range("L25").Value = 9.567733037336074
$ws.Range("O25").Value = 37.17113856588036
